$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 100005280
$ws.Range("I64").Value = 5147.5713
$ws.Range("J64").Value = 333338940
$ws.Range("K64").Value = 5147.5713
$ws.Range("L64").Value = 333338940
$ws.Range("M64").Value = -4899.5713
$ws.Range("N64").Value = -333339436

$ws.Range("H67").Value = 100005280
$ws.Range("I67").Value = 5147.5713
$ws.Range("J67").Value = 333338940
$ws.Range("K67").Value = 5147.5713
$ws.Range("L67").Value = 333338940
$ws.Range("M67").Value = -4289.5713
$ws.Range("N67").Value = -333340656

$ws.Range("H74").Value = 14157.315
$ws.Range("I74").Value = 15480.625
$ws.Range("J74").Value = 7099.6665
$ws.Range("K74").Value = 15480.625
$ws.Range("L74").Value = 7099.6665
$ws.Range("M74").Value = -14544.625
$ws.Range("N74").Value = -8971.666499999999

$ws.Range("H77").Value = 14157.315
$ws.Range("I77").Value = 15480.625
$ws.Range("J77").Value = 7099.6665
$ws.Range("K77").Value = 77403.125
$ws.Range("L77").Value = 35498.3325
$ws.Range("M77").Value = -72723.125
$ws.Range("N77").Value = -44858.3325

$ws.Range("H88").Value = 5886129.5
$ws.Range("I88").Value = 33334616
$ws.Range("J88").Value = 4311.0713
$ws.Range("K88").Value = 33334616
$ws.Range("L88").Value = 4311.0713
$ws.Range("M88").Value = -33334210
$ws.Range("N88").Value = -5123.0713

$ws.Range("H91").Value = 5886129.5
$ws.Range("I91").Value = 33334616
$ws.Range("J91").Value = 4311.0713
$ws.Range("K91").Value = 33334616
$ws.Range("L91").Value = 4311.0713
$ws.Range("M91").Value = -33333212
$ws.Range("N91").Value = -7119.0713

$ws.Range("H100").Value = 1449.75
$ws.Range("I100").Value = 933
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 933
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -392
$ws.Range("N100").Value = -4082

$ws.Range("H116").Value = 7966.4
$ws.Range("I116").Value = 4671.2856
$ws.Range("J116").Value = 10849.625
$ws.Range("K116").Value = 4671.2856
$ws.Range("L116").Value = 10849.625
$ws.Range("M116").Value = -1229.2856
$ws.Range("N116").Value = -17733.625

$ws.Range("H138").Value = 2895.258
$ws.Range("I138").Value = 2506.762
$ws.Range("J138").Value = 3094.244
$ws.Range("K138").Value = 7520.286
$ws.Range("L138").Value = 9282.732
$ws.Range("M138").Value = -2380.286
$ws.Range("N138").Value = -19562.732

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2829
$ws.Range("I26").Value = 993.5
$ws.Range("J26").Value = 6500
$ws.Range("K26").Value = 993.5
$ws.Range("L26").Value = 6500
$ws.Range("M26").Value = -663.5
$ws.Range("N26").Value = -7160

$ws.Range("H32").Value = 2883.121
$ws.Range("I32").Value = 2883.121
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2883.121
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2596.121
$ws.Range("N32").ClearContents()

$ws.Range("H74").Value = 223741.72
$ws.Range("I74").Value = 327971.2
$ws.Range("J74").Value = 2254.125
$ws.Range("K74").Value = 327971.2
$ws.Range("L74").Value = 2254.125
$ws.Range("M74").Value = -327097.2
$ws.Range("N74").Value = -4002.125

$ws.Range("H77").Value = 223741.72
$ws.Range("I77").Value = 327971.2
$ws.Range("J77").Value = 2254.125
$ws.Range("K77").Value = 1639856
$ws.Range("L77").Value = 11270.625
$ws.Range("M77").Value = -1635488
$ws.Range("N77").Value = -20006.625

$ws.Range("H102").Value = 2611.7727
$ws.Range("I102").Value = 2050.5789
$ws.Range("J102").Value = 6166
$ws.Range("K102").Value = 2050.5789
$ws.Range("L102").Value = 6166
$ws.Range("M102").Value = -428.5789
$ws.Range("N102").Value = -9410

$ws.Range("H132").Value = 2590.3333
$ws.Range("I132").Value = 2408.818
$ws.Range("J132").Value = 2953.3635
$ws.Range("K132").Value = 7226.454000000001
$ws.Range("L132").Value = 8860.0905
$ws.Range("M132").Value = -4696.454000000001
$ws.Range("N132").Value = -13920.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 33339832
$ws.Range("I20").Value = 62509924
$ws.Range("J20").Value = 2583.1428
$ws.Range("K20").Value = 62509924
$ws.Range("L20").Value = 2583.1428
$ws.Range("M20").Value = -62509677
$ws.Range("N20").Value = -3077.1428

$ws.Range("H64").Value = 1028.2858
$ws.Range("I64").Value = 626.5
$ws.Range("J64").Value = 1189
$ws.Range("K64").Value = 626.5
$ws.Range("L64").Value = 1189
$ws.Range("M64").Value = -401.5
$ws.Range("N64").Value = -1639

$ws.Range("H67").Value = 1028.2858
$ws.Range("I67").Value = 626.5
$ws.Range("J67").Value = 1189
$ws.Range("K67").Value = 626.5
$ws.Range("L67").Value = 1189
$ws.Range("M67").Value = 153.5
$ws.Range("N67").Value = -2749

$ws.Range("H105").Value = 20002000
$ws.Range("I105").Value = 1430013.9
$ws.Range("J105").Value = 41669316
$ws.Range("K105").Value = 1430013.9
$ws.Range("L105").Value = 41669316
$ws.Range("M105").Value = -1428266.9
$ws.Range("N105").Value = -41672810

$ws.Range("H134").Value = 2262.7346
$ws.Range("I134").Value = 1954.075
$ws.Range("J134").Value = 3634.5557
$ws.Range("K134").Value = 5862.225
$ws.Range("L134").Value = 10903.6671
$ws.Range("M134").Value = -3327.225
$ws.Range("N134").Value = -15973.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 36999
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 36999
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 36999
$ws.Range("N28").Value = -37489

$ws.Range("H58").Value = 1765.1111
$ws.Range("I58").Value = 1124.069
$ws.Range("J58").Value = 2927
$ws.Range("K58").Value = 1124.069
$ws.Range("L58").Value = 2927
$ws.Range("M58").Value = -921.069
$ws.Range("N58").Value = -3333

$ws.Range("H62").Value = 12503649
$ws.Range("I62").Value = 12503649
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 12503649
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -12503025

$ws.Range("H65").Value = 12503649
$ws.Range("I65").Value = 12503649
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 62518245
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -62515125

$ws.Range("H124").Value = 46250
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 46250
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 46250
$ws.Range("N124").Value = -51160

$ws.Range("H134").Value = 2413.5
$ws.Range("I134").Value = 2093.5933
$ws.Range("J134").Value = 4129.364
$ws.Range("K134").Value = 6280.7799
$ws.Range("L134").Value = 12388.092
$ws.Range("M134").Value = -3745.7799
$ws.Range("N134").Value = -17458.092

$ws.Range("H136").Value = 1765.1111
$ws.Range("I136").Value = 1124.069
$ws.Range("J136").Value = 2927
$ws.Range("K136").Value = 3372.207
$ws.Range("L136").Value = 8781
$ws.Range("M136").Value = -822.2069999999999
$ws.Range("N136").Value = -13881

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11088895
$ws.Range("I4").Value = 11144217
$ws.Range("J4").Value = 9816502
$ws.Range("K4").Value = 33432651
$ws.Range("L4").Value = 29449506
$ws.Range("M4").Value = -33432539
$ws.Range("N4").Value = -29449730

$ws.Range("H11").Value = 199
$ws.Range("I11").Value = 199
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 597
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -457
$ws.Range("N11").ClearContents()

$ws.Range("H51").Value = 3375
$ws.Range("I51").Value = 187.5
$ws.Range("J51").Value = 5500
$ws.Range("K51").Value = 562.5
$ws.Range("L51").Value = 16500
$ws.Range("M51").Value = -102.5
$ws.Range("N51").Value = -17420

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 27579.6
$ws.Range("I70").Value = 8999
$ws.Range("J70").Value = 32224.75
$ws.Range("K70").Value = 8999
$ws.Range("L70").Value = 32224.75
$ws.Range("M70").Value = -8729
$ws.Range("N70").Value = -32764.75

$ws.Range("H73").Value = 27579.6
$ws.Range("I73").Value = 8999
$ws.Range("J73").Value = 32224.75
$ws.Range("K73").Value = 8999
$ws.Range("L73").Value = 32224.75
$ws.Range("M73").Value = -8063
$ws.Range("N73").Value = -34096.75

$ws.Range("H132").Value = 2579.4119
$ws.Range("I132").Value = 2450
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 7350
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -4820
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4040.6316
$ws.Range("I7").Value = 2889.9092
$ws.Range("J7").Value = 5622.875
$ws.Range("K7").Value = 2889.9092
$ws.Range("L7").Value = 5622.875
$ws.Range("M7").Value = -2777.9092
$ws.Range("N7").Value = -5846.875

$ws.Range("H93").Value = 1637.2142
$ws.Range("I93").Value = 1592.5
$ws.Range("J93").Value = 1749
$ws.Range("K93").Value = 1592.5
$ws.Range("L93").Value = 1749
$ws.Range("M93").Value = -344.5
$ws.Range("N93").Value = -4245

$ws.Range("H126").Value = 4040.6316
$ws.Range("I126").Value = 2889.9092
$ws.Range("J126").Value = 5622.875
$ws.Range("K126").Value = 8669.7276
$ws.Range("L126").Value = 16868.625
$ws.Range("M126").Value = -6199.7276
$ws.Range("N126").Value = -21808.625

$ws.Range("H132").Value = 5879.136
$ws.Range("I132").Value = 2373.7273
$ws.Range("J132").Value = 9384.546
$ws.Range("K132").Value = 7121.1819
$ws.Range("L132").Value = 28153.638
$ws.Range("M132").Value = -4591.1819
$ws.Range("N132").Value = -33213.638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 38222.11
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 38222.11
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 38222.11
$ws.Range("N46").Value = -38684.11

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H96").Value = 3824.8
$ws.Range("I96").Value = 2321.2856
$ws.Range("J96").Value = 7333
$ws.Range("K96").Value = 2321.2856
$ws.Range("L96").Value = 7333
$ws.Range("M96").Value = -948.2856000000002
$ws.Range("N96").Value = -10079

$ws.Range("H122").Value = 35715680
$ws.Range("I122").Value = 1636
$ws.Range("J122").Value = 83334410
$ws.Range("K122").Value = 4908
$ws.Range("L122").Value = 250003230
$ws.Range("M122").Value = -2458
$ws.Range("N122").Value = -250008130

$ws.Range("H132").Value = 3156.9048
$ws.Range("I132").Value = 3331.3157
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 9993.947100000001
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -7463.947100000001
$ws.Range("N132").Value = -9560

$ws.Range("H134").Value = 38222.11
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 38222.11
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 114666.33
$ws.Range("N134").Value = -119736.33

